# Agregar formato de moneda a los totales
# - Insert a new "PLANILLA" header column (P2), copying O2's formatting
# - Rename header "TOTAL" (O2) to "TOTAL P."
# - Update the sheet view (top-left cell / active selection)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone O2 (value + formatting) into the new P2 column so it matches the
# existing header look & feel exactly.
$ws.Range("O2").Copy($ws.Range("P2"))

# Update header labels (PLANILLA first, then TOTAL P., to mirror the order
# the new shared strings were recorded in).
$ws.Range("P2").Value = "PLANILLA"
$ws.Range("O2").Value = "TOTAL P."

# Move the viewport / selection the way the saved workbook shows it.
$excel.ActiveWindow.ScrollColumn = 10
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("O9").Select()
